$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "cedula" column header and data
$ws.Range("E1").Value = "cedula"
$ws.Range("E2").Value = 123123123
$ws.Range("E3").Value = 12312313

# Row 4 has an empty, underlined cell (matches style added by author)
$ws.Range("E4").Font.Underline = $true

# Update selection to reflect where the user left off editing
$ws.Range("G5:H5").Select() | Out-Null

# Configure page setup (paper size / orientation)
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
